# Update "想去人数" (F column) and G18 value on both the "展览" and "全部类型"
# worksheets, which carry duplicate data tables.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 3186
    $ws.Range("F7").Value = 1684

    # G18 changes from a numeric price (50) to the text "不可售" (not for sale)
    $ws.Range("G18").Value = "不可售"

    $ws.Range("F19").Value = 25
    $ws.Range("F23").Value = 378
    $ws.Range("F24").Value = 209
    $ws.Range("F25").Value = 104
    $ws.Range("F29").Value = 243
    $ws.Range("F30").Value = 2161
    $ws.Range("F33").Value = 469
    $ws.Range("F34").Value = 326
    $ws.Range("F38").Value = 345
    $ws.Range("F40").Value = 516
}
